$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(98, 8).Value = 2000
$ws.Cells.Item(98, 9).Value = 2000
$ws.Cells.Item(98, 10).Value = 2000
$ws.Cells.Item(98, 11).Value = 2000
$ws.Cells.Item(98, 12).Value = 2000
$ws.Cells.Item(98, 13).Value = -502
$ws.Cells.Item(98, 14).Value = -4996
$ws.Cells.Item(106, 8).Value = 133337690
$ws.Cells.Item(106, 9).Value = 37041870
$ws.Cells.Item(106, 10).Value = 1000000000
$ws.Cells.Item(106, 11).Value = 37041870
$ws.Cells.Item(106, 12).Value = 1000000000
$ws.Cells.Item(106, 13).Value = -37041239
$ws.Cells.Item(106, 14).Value = -1000001262
$ws.Cells.Item(107, 8).Value = 17861898
$ws.Cells.Item(107, 9).Value = 22728590
$ws.Cells.Item(107, 11).Value = 22728590
$ws.Cells.Item(107, 13).Value = -22726670
$ws.Cells.Item(122, 8).Value = 2000
$ws.Cells.Item(122, 9).Value = 2000
$ws.Cells.Item(122, 10).Value = 2000
$ws.Cells.Item(122, 11).Value = 6000
$ws.Cells.Item(122, 12).Value = 6000
$ws.Cells.Item(122, 13).Value = -3550
$ws.Cells.Item(122, 14).Value = -10900
$ws.Cells.Item(129, 8).Value = 1042.8873
$ws.Cells.Item(129, 9).Value = 777.5
$ws.Cells.Item(129, 10).Value = 1058.7313
$ws.Cells.Item(129, 11).Value = 2332.5
$ws.Cells.Item(129, 12).Value = 3176.1939
$ws.Cells.Item(129, 13).Value = 2667.5
$ws.Cells.Item(129, 14).Value = -13176.1939
$ws.Cells.Item(135, 8).Value = 2812.3
$ws.Cells.Item(135, 9).Value = 2749.7896
$ws.Cells.Item(135, 10).Value = 4000
$ws.Cells.Item(135, 11).Value = 24748.1064
$ws.Cells.Item(135, 12).Value = 36000
$ws.Cells.Item(135, 13).Value = -22213.1064
$ws.Cells.Item(135, 14).Value = -41070
$ws.Cells.Item(138, 8).Value = 3698.8708
$ws.Cells.Item(138, 9).Value = 1857.8422
$ws.Cells.Item(138, 10).Value = 4512.3486
$ws.Cells.Item(138, 11).Value = 5573.5266
$ws.Cells.Item(138, 12).Value = 13537.0458
$ws.Cells.Item(138, 13).Value = -433.5266000000001
$ws.Cells.Item(138, 14).Value = -23817.0458

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 1431.45
$ws.Cells.Item(2, 9).Value = 1455.2106
$ws.Cells.Item(2, 11).Value = 1455.2106
$ws.Cells.Item(2, 13).Value = -1342.2106
$ws.Cells.Item(32, 8).Value = 3570.5657
$ws.Cells.Item(32, 9).Value = 3083.5964
$ws.Cells.Item(32, 10).Value = 5031.4736
$ws.Cells.Item(32, 11).Value = 3083.5964
$ws.Cells.Item(32, 12).Value = 5031.4736
$ws.Cells.Item(32, 13).Value = -2796.5964
$ws.Cells.Item(32, 14).Value = -5605.4736
$ws.Cells.Item(45, 8).Value = 9399.556
$ws.Cells.Item(45, 9).Value = 11645.571
$ws.Cells.Item(45, 10).Value = 1538.5
$ws.Cells.Item(45, 11).Value = 11645.571
$ws.Cells.Item(45, 12).Value = 1538.5
$ws.Cells.Item(45, 13).Value = -11268.571
$ws.Cells.Item(45, 14).Value = -2292.5
$ws.Cells.Item(110, 8).Value = 1261.5
$ws.Cells.Item(110, 9).Value = 1245.7778
$ws.Cells.Item(110, 11).Value = 1245.7778
$ws.Cells.Item(110, 13).Value = 799.2221999999999
$ws.Cells.Item(116, 8).Value = 1431.45
$ws.Cells.Item(116, 9).Value = 1455.2106
$ws.Cells.Item(116, 11).Value = 1455.2106
$ws.Cells.Item(116, 13).Value = 838.7893999999999
$ws.Cells.Item(122, 8).Value = 1711226
$ws.Cells.Item(122, 9).Value = 2850321
$ws.Cells.Item(122, 10).Value = 2583.3333
$ws.Cells.Item(122, 11).Value = 8550963
$ws.Cells.Item(122, 12).Value = 7749.999899999999
$ws.Cells.Item(122, 13).Value = -8548513
$ws.Cells.Item(122, 14).Value = -12649.9999
$ws.Cells.Item(132, 8).Value = 1886.3818
$ws.Cells.Item(132, 9).Value = 1607.2554
$ws.Cells.Item(132, 10).Value = 3526.25
$ws.Cells.Item(132, 11).Value = 4821.7662
$ws.Cells.Item(132, 12).Value = 10578.75
$ws.Cells.Item(132, 13).Value = -2291.7662
$ws.Cells.Item(132, 14).Value = -15638.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 1431.45
$ws.Cells.Item(3, 9).Value = 1455.2106
$ws.Cells.Item(3, 11).Value = 1455.2106
$ws.Cells.Item(3, 13).Value = -1341.2106
$ws.Cells.Item(47, 8).Value = 150000
$ws.Cells.Item(47, 10).Value = 150000
$ws.Cells.Item(47, 12).Value = 150000
$ws.Cells.Item(47, 14).Value = -151040
$ws.Cells.Item(107, 8).Value = 1159
$ws.Cells.Item(107, 9).Value = 1148.6666
$ws.Cells.Item(107, 10).Value = 1190
$ws.Cells.Item(107, 11).Value = 1148.6666
$ws.Cells.Item(107, 12).Value = 1190
$ws.Cells.Item(107, 13).Value = 771.3334
$ws.Cells.Item(107, 14).Value = -5030

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 4512.3335
$ws.Cells.Item(16, 9).Value = 991
$ws.Cells.Item(16, 10).Value = 6273
$ws.Cells.Item(16, 11).Value = 991
$ws.Cells.Item(16, 12).Value = 6273
$ws.Cells.Item(16, 13).Value = -704
$ws.Cells.Item(16, 14).Value = -6847
$ws.Cells.Item(113, 8).Value = 4512.3335
$ws.Cells.Item(113, 9).Value = 991
$ws.Cells.Item(113, 10).Value = 6273
$ws.Cells.Item(113, 11).Value = 991
$ws.Cells.Item(113, 12).Value = 6273
$ws.Cells.Item(113, 13).Value = 1179
$ws.Cells.Item(113, 14).Value = -10613
$ws.Cells.Item(122, 8).Value = 1769.1765
$ws.Cells.Item(122, 9).Value = 1286.6154
$ws.Cells.Item(122, 10).Value = 3337.5
$ws.Cells.Item(122, 11).Value = 3859.8462
$ws.Cells.Item(122, 12).Value = 10012.5
$ws.Cells.Item(122, 13).Value = -1409.8462
$ws.Cells.Item(122, 14).Value = -14912.5
$ws.Cells.Item(132, 8).Value = 2216.0334
$ws.Cells.Item(132, 9).Value = 1899.3462
$ws.Cells.Item(132, 10).Value = 4274.5
$ws.Cells.Item(132, 11).Value = 5698.0386
$ws.Cells.Item(132, 12).Value = 12823.5
$ws.Cells.Item(132, 13).Value = -3168.0386
$ws.Cells.Item(132, 14).Value = -17883.5
$ws.Cells.Item(134, 8).Value = 3307.2
$ws.Cells.Item(134, 9).Value = 3307.2
$ws.Cells.Item(134, 11).Value = 9921.599999999999
$ws.Cells.Item(134, 13).Value = -7386.599999999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value = 4082.75
$ws.Cells.Item(68, 9).Value = 8090.4287
$ws.Cells.Item(68, 10).Value = 1532.409
$ws.Cells.Item(68, 11).Value = 24271.2861
$ws.Cells.Item(68, 12).Value = 4597.227000000001
$ws.Cells.Item(68, 13).Value = -23460.2861
$ws.Cells.Item(68, 14).Value = -6219.227000000001
$ws.Cells.Item(71, 8).Value = 4082.75
$ws.Cells.Item(71, 9).Value = 8090.4287
$ws.Cells.Item(71, 10).Value = 1532.409
$ws.Cells.Item(71, 11).Value = 72813.85830000001
$ws.Cells.Item(71, 12).Value = 13791.681
$ws.Cells.Item(71, 13).Value = -68757.85830000001
$ws.Cells.Item(71, 14).Value = -21903.681
$ws.Cells.Item(122, 8).Value = 600
$ws.Cells.Item(122, 9).Value = 400
$ws.Cells.Item(122, 10).Value = 800
$ws.Cells.Item(122, 11).Value = 3600
$ws.Cells.Item(122, 12).Value = 7200
$ws.Cells.Item(122, 13).Value = -1150
$ws.Cells.Item(122, 14).Value = -12100

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5605.5645
$ws.Cells.Item(70, 9).Value = 5648.9585
$ws.Cells.Item(70, 10).Value = 5456.7856
$ws.Cells.Item(70, 11).Value = 5648.9585
$ws.Cells.Item(70, 12).Value = 5456.7856
$ws.Cells.Item(70, 13).Value = -5378.9585
$ws.Cells.Item(70, 14).Value = -5996.7856
$ws.Cells.Item(73, 8).Value = 5605.5645
$ws.Cells.Item(73, 9).Value = 5648.9585
$ws.Cells.Item(73, 10).Value = 5456.7856
$ws.Cells.Item(73, 11).Value = 5648.9585
$ws.Cells.Item(73, 12).Value = 5456.7856
$ws.Cells.Item(73, 13).Value = -4712.9585
$ws.Cells.Item(73, 14).Value = -7328.7856
$ws.Cells.Item(102, 8).Value = 841.8570999999999
$ws.Cells.Item(102, 9).Value = 515.5
$ws.Cells.Item(102, 11).Value = 515.5
$ws.Cells.Item(102, 13).Value = 1106.5
$ws.Cells.Item(113, 8).Value = 55556484
$ws.Cells.Item(113, 9).Value = 66667508
$ws.Cells.Item(113, 11).Value = 66667508
$ws.Cells.Item(113, 13).Value = -66665338
$ws.Cells.Item(122, 8).Value = 152117820
$ws.Cells.Item(122, 9).Value = 212963760
$ws.Cells.Item(122, 10).Value = 3000
$ws.Cells.Item(122, 11).Value = 638891280
$ws.Cells.Item(122, 12).Value = 9000
$ws.Cells.Item(122, 13).Value = -638888830
$ws.Cells.Item(122, 14).Value = -13900
$ws.Cells.Item(132, 8).Value = 3303.1633
$ws.Cells.Item(132, 9).Value = 2882.484
$ws.Cells.Item(132, 10).Value = 4027.6667
$ws.Cells.Item(132, 11).Value = 8647.451999999999
$ws.Cells.Item(132, 12).Value = 12083.0001
$ws.Cells.Item(132, 13).Value = -6117.451999999999
$ws.Cells.Item(132, 14).Value = -17143.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2857.7144
$ws.Cells.Item(7, 9).Value = 2168
$ws.Cells.Item(7, 10).Value = 3375
$ws.Cells.Item(7, 11).Value = 2168
$ws.Cells.Item(7, 12).Value = 3375
$ws.Cells.Item(7, 13).Value = -2056
$ws.Cells.Item(7, 14).Value = -3599
$ws.Cells.Item(40, 8).Value = 1000000000
$ws.Cells.Item(40, 9).Value = 1000000000
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 11).Value = 1000000000
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 13).Value = -999999864
$ws.Cells.Item(40, 14).Value = ""
$ws.Cells.Item(61, 8).Value = 2530.6875
$ws.Cells.Item(61, 9).Value = 2535.0715
$ws.Cells.Item(61, 10).Value = 2500
$ws.Cells.Item(61, 11).Value = 2535.0715
$ws.Cells.Item(61, 12).Value = 2500
$ws.Cells.Item(61, 13).Value = -2333.0715
$ws.Cells.Item(61, 14).Value = -2904
$ws.Cells.Item(113, 8).Value = 2530.6875
$ws.Cells.Item(113, 9).Value = 2535.0715
$ws.Cells.Item(113, 10).Value = 2500
$ws.Cells.Item(113, 11).Value = 2535.0715
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 13).Value = -365.0715
$ws.Cells.Item(113, 14).Value = -6840
$ws.Cells.Item(122, 8).Value = 5430889.5
$ws.Cells.Item(122, 9).Value = 5496795.5
$ws.Cells.Item(122, 11).Value = 16490386.5
$ws.Cells.Item(122, 13).Value = -16487936.5
$ws.Cells.Item(126, 8).Value = 2857.7144
$ws.Cells.Item(126, 9).Value = 2168
$ws.Cells.Item(126, 10).Value = 3375
$ws.Cells.Item(126, 11).Value = 6504
$ws.Cells.Item(126, 12).Value = 10125
$ws.Cells.Item(126, 13).Value = -4034
$ws.Cells.Item(126, 14).Value = -15065
$ws.Cells.Item(136, 8).Value = 3799.7954
$ws.Cells.Item(136, 10).Value = 6999.5454
$ws.Cells.Item(136, 12).Value = 20998.6362
$ws.Cells.Item(136, 14).Value = -26098.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 52632136
$ws.Cells.Item(107, 9).Value = 90909576
$ws.Cells.Item(107, 10).Value = 655.75
$ws.Cells.Item(107, 11).Value = 272728728
$ws.Cells.Item(107, 12).Value = 1967.25
$ws.Cells.Item(107, 13).Value = -272726808
$ws.Cells.Item(107, 14).Value = -5807.25
$ws.Cells.Item(113, 8).Value = 1753.7693
$ws.Cells.Item(113, 9).Value = 1566.5
$ws.Cells.Item(113, 11).Value = 4699.5
$ws.Cells.Item(113, 13).Value = -2529.5
$ws.Cells.Item(122, 8).Value = 2454.889
$ws.Cells.Item(122, 9).Value = 2136.75
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 6410.25
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -3960.25
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(126, 8).Value = 1026.75
$ws.Cells.Item(126, 9).Value = 1001
$ws.Cells.Item(126, 10).Value = 1052.5
$ws.Cells.Item(126, 11).Value = 3003
$ws.Cells.Item(126, 12).Value = 3157.5
$ws.Cells.Item(126, 13).Value = -533
$ws.Cells.Item(126, 14).Value = -8097.5
$ws.Cells.Item(132, 8).Value = 2222.7368
$ws.Cells.Item(132, 9).Value = 2019.0416
$ws.Cells.Item(132, 10).Value = 2571.9285
$ws.Cells.Item(132, 11).Value = 6057.1248
$ws.Cells.Item(132, 12).Value = 7715.7855
$ws.Cells.Item(132, 13).Value = -3527.1248
$ws.Cells.Item(132, 14).Value = -12775.7855
$ws.Cells.Item(136, 8).Value = 1010.3889
$ws.Cells.Item(136, 9).Value = 934.5294
$ws.Cells.Item(136, 10).Value = 2300
$ws.Cells.Item(136, 11).Value = 2803.5882
$ws.Cells.Item(136, 12).Value = 6900
$ws.Cells.Item(136, 13).Value = -253.5882000000001
$ws.Cells.Item(136, 14).Value = -12000
